# Rombak sistem validasi file pake web socket
#
# 1) "MERR" -> "__DAERAHBANK__" (single run's text, inside the assignment
#    paragraph that also carries __PENGGANTIKANDIDAT__ / __TANGGALMULAI__ /
#    __TANGGALSELESAI__ placeholders).
# 2) "Ma’aliy" split into its own run bracketed by spellStart/spellEnd
#    proofErr markers (leaving the leading space as a separate, unmarked run),
#    matching the signature-block paragraph "M. Dlou’ul Ma’aliy".
#
# Note on this runtime: Range.InsertXML() replaces the *entire* paragraph
# that the target Range lives in (there is no way to splice a bare <w:r>/
# <w:proofErr> into the middle of an existing paragraph's run list), so each
# payload below is the complete, paragraph XML (pPr + every run) with only
# the minimal required change applied; every other run/attribute is
# byte-for-byte identical to what was already in the document.

$d = $word.ActiveDocument

function Get-ContainingParagraphRange($rng) {
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($rng.Start -ge $p.Range.Start -and $rng.End -le $p.Range.End) {
            return $p.Range
        }
    }
    throw "Get-ContainingParagraphRange: no containing paragraph found"
}

# --- Edit 1: MERR -> __DAERAHBANK__ --------------------------------------
$rng1 = $d.Content
$ok1 = $rng1.Find.Execute("MERR", $true, $false, $false, $false, $false,
                           $true, 1, $false, "", 0)
if (-not $ok1) { throw "MERR run not found" }
$para1 = Get-ContainingParagraphRange $rng1
$para1.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5ECC1AF2" w14:textId="2E412882" w:rsidR="00E90B7A" w:rsidRDefault="00E90B7A" w:rsidP="00E90B7A"><w:pPr><w:spacing w:after="0" w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Menugaskan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sebagai</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pengganti</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>Driver</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>atas nama</w:t></w:r><w:r w:rsidR="004D489D"><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve"> __PENGGANTIKANDIDAT__</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00261993"><w:t xml:space="preserve">di </w:t></w:r><w:r><w:t xml:space="preserve">PT. Bank </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Mandiri</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (Persero) </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Tbk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="004D489D"><w:t>__DAERAHBANK__</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve">pada tanggal </w:t></w:r><w:r w:rsidR="004D489D"><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>__TANGGALMULAI__</w:t></w:r><w:r w:rsidR="005B4D57"><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve"> sampai </w:t></w:r><w:r w:rsidR="004D489D"><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>__TANGGALSELESAI__</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# --- Edit 2: mark "Ma’aliy" with spellStart/spellEnd proofErr -------------
$rng2 = $d.Content
$ok2 = $rng2.Find.Execute("Ma’aliy", $true, $false, $false, $false, $false,
                           $true, 1, $false, "", 0)
if (-not $ok2) { throw "Ma’aliy run not found" }
$para2 = Get-ContainingParagraphRange $rng2
$para2.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="50772BDE" w14:textId="1B0764ED" w:rsidR="007B5FC3" w:rsidRDefault="007B5FC3" w:rsidP="00AA566A"><w:pPr><w:jc w:val="both"/><w:rPr><w:iCs/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">M. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Dlou’ul</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Ma’aliy</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null
